$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Tnc"
$ws.Cells.Item(2,3).Value = "Itga7"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value2 = 3
$ws.Cells.Item(2,6).Value2 = 1
$ws.Cells.Item(2,7).Value2 = 8.948174
$ws.Cells.Item(2,8).Value2 = 26.844522
$ws.Cells.Item(2,9).Value2 = 0.0695931738232498
$ws.Cells.Item(2,10).Value2 = 0.0695931738232498
$ws.Cells.Item(2,11).Value2 = 3
$ws.Cells.Item(2,12).Value2 = 1
$ws.Cells.Item(2,13).Value2 = 5.636552333333334
$ws.Cells.Item(2,14).Value2 = 16.909657
$ws.Cells.Item(2,15).Value2 = 0.09690140221006956
$ws.Cells.Item(2,16).Value2 = 0.09690140221006956
$ws.Cells.Item(2,17).Value2 = 50.43685103877267
$ws.Cells.Item(2,18).Value2 = 453.931659348954
$ws.Cells.Item(2,19).Value2 = 0.006743676127722012
$ws.Cells.Item(2,20).Value2 = 0.006743676127722012

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Tnc"
$ws.Cells.Item(3,3).Value = "Itga7"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value2 = 3
$ws.Cells.Item(3,6).Value2 = 1
$ws.Cells.Item(3,7).Value2 = 8.948174
$ws.Cells.Item(3,8).Value2 = 26.844522
$ws.Cells.Item(3,9).Value2 = 0.0695931738232498
$ws.Cells.Item(3,10).Value2 = 0.0695931738232498
$ws.Cells.Item(3,11).Value2 = 3
$ws.Cells.Item(3,12).Value2 = 1
$ws.Cells.Item(3,13).Value2 = 4.050135999999999
$ws.Cells.Item(3,14).Value2 = 12.150408
$ws.Cells.Item(3,15).Value2 = 0.06962835335006774
$ws.Cells.Item(3,16).Value2 = 0.06962835335006774
$ws.Cells.Item(3,17).Value2 = 36.24132165166399
$ws.Cells.Item(3,18).Value2 = 326.171894864976
$ws.Cells.Item(3,19).Value2 = 0.004845658097717922
$ws.Cells.Item(3,20).Value2 = 0.004845658097717922

$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Tnc"
$ws.Cells.Item(4,3).Value = "Itga7"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value2 = 3
$ws.Cells.Item(4,6).Value2 = 1
$ws.Cells.Item(4,7).Value2 = 8.948174
$ws.Cells.Item(4,8).Value2 = 26.844522
$ws.Cells.Item(4,9).Value2 = 0.0695931738232498
$ws.Cells.Item(4,10).Value2 = 0.0695931738232498
$ws.Cells.Item(4,11).Value2 = 3
$ws.Cells.Item(4,12).Value2 = 1
$ws.Cells.Item(4,13).Value2 = 0.6741636666666667
$ws.Cells.Item(4,14).Value2 = 2.022491
$ws.Cells.Item(4,15).Value2 = 0.01158995796645939
$ws.Cells.Item(4,16).Value2 = 0.01158995796645939
$ws.Cells.Item(4,17).Value2 = 6.032533793811334
$ws.Cells.Item(4,18).Value2 = 54.292804144302
$ws.Cells.Item(4,19).Value2 = 0.0008065819593639671
$ws.Cells.Item(4,20).Value2 = 0.0008065819593639669

$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Tnc"
$ws.Cells.Item(5,3).Value = "Itga7"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value2 = 3
$ws.Cells.Item(5,6).Value2 = 1
$ws.Cells.Item(5,7).Value2 = 8.948174
$ws.Cells.Item(5,8).Value2 = 26.844522
$ws.Cells.Item(5,9).Value2 = 0.0695931738232498
$ws.Cells.Item(5,10).Value2 = 0.0695931738232498
$ws.Cells.Item(5,11).Value2 = 3
$ws.Cells.Item(5,12).Value2 = 1
$ws.Cells.Item(5,13).Value2 = 47.807061
$ws.Cells.Item(5,14).Value2 = 143.421183
$ws.Cells.Item(5,15).Value2 = 0.8218802864734033
$ws.Cells.Item(5,16).Value2 = 0.8218802864734033
$ws.Cells.Item(5,17).Value2 = 427.785900256614
$ws.Cells.Item(5,18).Value2 = 3850.073102309526
$ws.Cells.Item(5,19).Value2 = 0.05719725763844589
$ws.Cells.Item(5,20).Value2 = 0.05719725763844589

$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Tnc"
$ws.Cells.Item(6,3).Value = "Itga7"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value2 = 3
$ws.Cells.Item(6,6).Value2 = 1
$ws.Cells.Item(6,7).Value2 = 105.8801323333333
$ws.Cells.Item(6,8).Value2 = 317.640397
$ws.Cells.Item(6,9).Value2 = 0.8234679448457706
$ws.Cells.Item(6,10).Value2 = 0.8234679448457706
$ws.Cells.Item(6,11).Value2 = 3
$ws.Cells.Item(6,12).Value2 = 1
$ws.Cells.Item(6,13).Value2 = 5.636552333333334
$ws.Cells.Item(6,14).Value2 = 16.909657
$ws.Cells.Item(6,15).Value2 = 0.09690140221006956
$ws.Cells.Item(6,16).Value2 = 0.09690140221006956
$ws.Cells.Item(6,17).Value2 = 596.7989069570922
$ws.Cells.Item(6,18).Value2 = 5371.19016261383
$ws.Cells.Item(6,19).Value2 = 0.07979519853059938
$ws.Cells.Item(6,20).Value2 = 0.07979519853059938

$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Tnc"
$ws.Cells.Item(7,3).Value = "Itga7"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value2 = 3
$ws.Cells.Item(7,6).Value2 = 1
$ws.Cells.Item(7,7).Value2 = 105.8801323333333
$ws.Cells.Item(7,8).Value2 = 317.640397
$ws.Cells.Item(7,9).Value2 = 0.8234679448457706
$ws.Cells.Item(7,10).Value2 = 0.8234679448457706
$ws.Cells.Item(7,11).Value2 = 3
$ws.Cells.Item(7,12).Value2 = 1
$ws.Cells.Item(7,13).Value2 = 4.050135999999999
$ws.Cells.Item(7,14).Value2 = 12.150408
$ws.Cells.Item(7,15).Value2 = 0.06962835335006774
$ws.Cells.Item(7,16).Value2 = 0.06962835335006774
$ws.Cells.Item(7,17).Value2 = 428.8289356479972
$ws.Cells.Item(7,18).Value2 = 3859.460420831976
$ws.Cells.Item(7,19).Value2 = 0.05733671703617541
$ws.Cells.Item(7,20).Value2 = 0.05733671703617541

$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Tnc"
$ws.Cells.Item(8,3).Value = "Itga7"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(8,5).Value2 = 3
$ws.Cells.Item(8,6).Value2 = 1
$ws.Cells.Item(8,7).Value2 = 105.8801323333333
$ws.Cells.Item(8,8).Value2 = 317.640397
$ws.Cells.Item(8,9).Value2 = 0.8234679448457706
$ws.Cells.Item(8,10).Value2 = 0.8234679448457706
$ws.Cells.Item(8,11).Value2 = 3
$ws.Cells.Item(8,12).Value2 = 1
$ws.Cells.Item(8,13).Value2 = 0.6741636666666667
$ws.Cells.Item(8,14).Value2 = 2.022491
$ws.Cells.Item(8,15).Value2 = 0.01158995796645939
$ws.Cells.Item(8,16).Value2 = 0.01158995796645939
$ws.Cells.Item(8,17).Value2 = 71.3805382409919
$ws.Cells.Item(8,18).Value2 = 642.424844168927
$ws.Cells.Item(8,19).Value2 = 0.00954395886748918
$ws.Cells.Item(8,20).Value2 = 0.009543958867489178

$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Tnc"
$ws.Cells.Item(9,3).Value = "Itga7"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value2 = 3
$ws.Cells.Item(9,6).Value2 = 1
$ws.Cells.Item(9,7).Value2 = 105.8801323333333
$ws.Cells.Item(9,8).Value2 = 317.640397
$ws.Cells.Item(9,9).Value2 = 0.8234679448457706
$ws.Cells.Item(9,10).Value2 = 0.8234679448457706
$ws.Cells.Item(9,11).Value2 = 3
$ws.Cells.Item(9,12).Value2 = 1
$ws.Cells.Item(9,13).Value2 = 47.807061
$ws.Cells.Item(9,14).Value2 = 143.421183
$ws.Cells.Item(9,15).Value2 = 0.8218802864734033
$ws.Cells.Item(9,16).Value2 = 0.8218802864734033
$ws.Cells.Item(9,17).Value2 = 5061.81794514774
$ws.Cells.Item(9,18).Value2 = 45556.36150632965
$ws.Cells.Item(9,19).Value2 = 0.6767920704115066
$ws.Cells.Item(9,20).Value2 = 0.6767920704115066

$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Tnc"
$ws.Cells.Item(10,3).Value = "Itga7"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value2 = 2
$ws.Cells.Item(10,6).Value2 = 0.6666666666666666
$ws.Cells.Item(10,7).Value2 = 0.1148696666666667
$ws.Cells.Item(10,8).Value2 = 0.344609
$ws.Cells.Item(10,9).Value2 = 0.0008933827928862465
$ws.Cells.Item(10,10).Value2 = 0.0008933827928862465
$ws.Cells.Item(10,11).Value2 = 3
$ws.Cells.Item(10,12).Value2 = 1
$ws.Cells.Item(10,13).Value2 = 5.636552333333334
$ws.Cells.Item(10,14).Value2 = 16.909657
$ws.Cells.Item(10,15).Value2 = 0.09690140221006956
$ws.Cells.Item(10,16).Value2 = 0.09690140221006956
$ws.Cells.Item(10,17).Value2 = 0.6474688876792223
$ws.Cells.Item(10,18).Value2 = 5.827219989113001
$ws.Cells.Item(10,19).Value2 = 0.00008657004534102544
$ws.Cells.Item(10,20).Value2 = 0.00008657004534102544

$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,2).Value = "Tnc"
$ws.Cells.Item(11,3).Value = "Itga7"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value2 = 2
$ws.Cells.Item(11,6).Value2 = 0.6666666666666666
$ws.Cells.Item(11,7).Value2 = 0.1148696666666667
$ws.Cells.Item(11,8).Value2 = 0.344609
$ws.Cells.Item(11,9).Value2 = 0.0008933827928862465
$ws.Cells.Item(11,10).Value2 = 0.0008933827928862465
$ws.Cells.Item(11,11).Value2 = 3
$ws.Cells.Item(11,12).Value2 = 1
$ws.Cells.Item(11,13).Value2 = 4.050135999999999
$ws.Cells.Item(11,14).Value2 = 12.150408
$ws.Cells.Item(11,15).Value2 = 0.06962835335006774
$ws.Cells.Item(11,16).Value2 = 0.06962835335006774
$ws.Cells.Item(11,17).Value2 = 0.4652377722746666
$ws.Cells.Item(11,18).Value2 = 4.187139950472
$ws.Cells.Item(11,19).Value2 = 0.00006220477277995396
$ws.Cells.Item(11,20).Value2 = 0.00006220477277995396

$ws.Cells.Item(12,1).Value = "M2"
$ws.Cells.Item(12,2).Value = "Tnc"
$ws.Cells.Item(12,3).Value = "Itga7"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value2 = 2
$ws.Cells.Item(12,6).Value2 = 0.6666666666666666
$ws.Cells.Item(12,7).Value2 = 0.1148696666666667
$ws.Cells.Item(12,8).Value2 = 0.344609
$ws.Cells.Item(12,9).Value2 = 0.0008933827928862465
$ws.Cells.Item(12,10).Value2 = 0.0008933827928862465
$ws.Cells.Item(12,11).Value2 = 3
$ws.Cells.Item(12,12).Value2 = 1
$ws.Cells.Item(12,13).Value2 = 0.6741636666666667
$ws.Cells.Item(12,14).Value2 = 2.022491
$ws.Cells.Item(12,15).Value2 = 0.01158995796645939
$ws.Cells.Item(12,16).Value2 = 0.01158995796645939
$ws.Cells.Item(12,17).Value2 = 0.07744095566877778
$ws.Cells.Item(12,18).Value2 = 0.696968601019
$ws.Cells.Item(12,19).Value2 = 0.00001035426901750969
$ws.Cells.Item(12,20).Value2 = 0.00001035426901750969

$ws.Cells.Item(13,1).Value = "M2"
$ws.Cells.Item(13,2).Value = "Tnc"
$ws.Cells.Item(13,3).Value = "Itga7"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value2 = 2
$ws.Cells.Item(13,6).Value2 = 0.6666666666666666
$ws.Cells.Item(13,7).Value2 = 0.1148696666666667
$ws.Cells.Item(13,8).Value2 = 0.344609
$ws.Cells.Item(13,9).Value2 = 0.0008933827928862465
$ws.Cells.Item(13,10).Value2 = 0.0008933827928862465
$ws.Cells.Item(13,11).Value2 = 3
$ws.Cells.Item(13,12).Value2 = 1
$ws.Cells.Item(13,13).Value2 = 47.807061
$ws.Cells.Item(13,14).Value2 = 143.421183
$ws.Cells.Item(13,15).Value2 = 0.8218802864734033
$ws.Cells.Item(13,16).Value2 = 0.8218802864734033
$ws.Cells.Item(13,17).Value2 = 5.491581161383
$ws.Cells.Item(13,18).Value2 = 49.42423045244701
$ws.Cells.Item(13,19).Value2 = 0.0007342537057477574
$ws.Cells.Item(13,20).Value2 = 0.0007342537057477574

$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "Tnc"
$ws.Cells.Item(14,3).Value = "Itga7"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value2 = 3
$ws.Cells.Item(14,6).Value2 = 1
$ws.Cells.Item(14,7).Value2 = 13.635153
$ws.Cells.Item(14,8).Value2 = 40.905459
$ws.Cells.Item(14,9).Value2 = 0.1060454985380935
$ws.Cells.Item(14,10).Value2 = 0.1060454985380935
$ws.Cells.Item(14,11).Value2 = 3
$ws.Cells.Item(14,12).Value2 = 1
$ws.Cells.Item(14,13).Value2 = 5.636552333333334
$ws.Cells.Item(14,14).Value2 = 16.909657
$ws.Cells.Item(14,15).Value2 = 0.09690140221006956
$ws.Cells.Item(14,16).Value2 = 0.09690140221006956
$ws.Cells.Item(14,17).Value2 = 76.85525345750702
$ws.Cells.Item(14,18).Value2 = 691.6972811175631
$ws.Cells.Item(14,19).Value2 = 0.01027595750640714
$ws.Cells.Item(14,20).Value2 = 0.01027595750640714

$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "Tnc"
$ws.Cells.Item(15,3).Value = "Itga7"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value2 = 3
$ws.Cells.Item(15,6).Value2 = 1
$ws.Cells.Item(15,7).Value2 = 13.635153
$ws.Cells.Item(15,8).Value2 = 40.905459
$ws.Cells.Item(15,9).Value2 = 0.1060454985380935
$ws.Cells.Item(15,10).Value2 = 0.1060454985380935
$ws.Cells.Item(15,11).Value2 = 3
$ws.Cells.Item(15,12).Value2 = 1
$ws.Cells.Item(15,13).Value2 = 4.050135999999999
$ws.Cells.Item(15,14).Value2 = 12.150408
$ws.Cells.Item(15,15).Value2 = 0.06962835335006774
$ws.Cells.Item(15,16).Value2 = 0.06962835335006774
$ws.Cells.Item(15,17).Value2 = 55.224224030808
$ws.Cells.Item(15,18).Value2 = 497.018016277272
$ws.Cells.Item(15,19).Value2 = 0.007383773443394465
$ws.Cells.Item(15,20).Value2 = 0.007383773443394464

$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "Tnc"
$ws.Cells.Item(16,3).Value = "Itga7"
$ws.Cells.Item(16,4).Value = "M2"
$ws.Cells.Item(16,5).Value2 = 3
$ws.Cells.Item(16,6).Value2 = 1
$ws.Cells.Item(16,7).Value2 = 13.635153
$ws.Cells.Item(16,8).Value2 = 40.905459
$ws.Cells.Item(16,9).Value2 = 0.1060454985380935
$ws.Cells.Item(16,10).Value2 = 0.1060454985380935
$ws.Cells.Item(16,11).Value2 = 3
$ws.Cells.Item(16,12).Value2 = 1
$ws.Cells.Item(16,13).Value2 = 0.6741636666666667
$ws.Cells.Item(16,14).Value2 = 2.022491
$ws.Cells.Item(16,15).Value2 = 0.01158995796645939
$ws.Cells.Item(16,16).Value2 = 0.01158995796645939
$ws.Cells.Item(16,17).Value2 = 9.192324742041
$ws.Cells.Item(16,18).Value2 = 82.730922678369
$ws.Cells.Item(16,19).Value2 = 0.001229062870588734
$ws.Cells.Item(16,20).Value2 = 0.001229062870588734

$ws.Cells.Item(17,1).Value = "sCs"
$ws.Cells.Item(17,2).Value = "Tnc"
$ws.Cells.Item(17,3).Value = "Itga7"
$ws.Cells.Item(17,4).Value = "sCs"
$ws.Cells.Item(17,5).Value2 = 3
$ws.Cells.Item(17,6).Value2 = 1
$ws.Cells.Item(17,7).Value2 = 13.635153
$ws.Cells.Item(17,8).Value2 = 40.905459
$ws.Cells.Item(17,9).Value2 = 0.1060454985380935
$ws.Cells.Item(17,10).Value2 = 0.1060454985380935
$ws.Cells.Item(17,11).Value2 = 3
$ws.Cells.Item(17,12).Value2 = 1
$ws.Cells.Item(17,13).Value2 = 47.807061
$ws.Cells.Item(17,14).Value2 = 143.421183
$ws.Cells.Item(17,15).Value2 = 0.8218802864734033
$ws.Cells.Item(17,16).Value2 = 0.8218802864734033
$ws.Cells.Item(17,17).Value2 = 651.8565912153331
$ws.Cells.Item(17,18).Value2 = 5866.709320937997
$ws.Cells.Item(17,19).Value2 = 0.08715670471770313
$ws.Cells.Item(17,20).Value2 = 0.08715670471770312
